# SCD0021 - Migrasi data Non Sales yang ada pada SAPM ke Digi Sales.xlsx
# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet to match the file's test-case id (SCD0313 -> SCD0021)
$ws.Name = "SCD0021"

# 2. Update the TC_ID column (B) rows 2-4 from the old "DGS-328" ticket id
#    to the new "SCD0021-001" test-case id.
$ws.Range("B2").Value = "SCD0021-001"
$ws.Range("B3").Value = "SCD0021-001"
$ws.Range("B4").Value = "SCD0021-001"

# 3. The longer TC_ID text no longer fits the old best-fit column width,
#    so widen column B to accommodate it.
$ws.Columns.Item(2).ColumnWidth = 11.65

# 4. Reset the view: scroll back to the left edge, zoom to 85%, and move
#    the active selection to B5.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B5").Select()
